# Update Sema5a-Met LR-pairs worksheet with new TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9750066666666667
$ws.Range("H2").Value = 2.92502
$ws.Range("I2").Value = 0.03314938328005726
$ws.Range("J2").Value = 0.03314938328005725
$ws.Range("M2").Value = 4.824089
$ws.Range("N2").Value = 14.472267
$ws.Range("O2").Value = 0.1181976021471384
$ws.Range("P2").Value = 0.1181976021471384
$ws.Range("Q2").Value = 4.703518935593333
$ws.Range("R2").Value = 42.33167042034
$ws.Range("S2").Value = 0.00391817761635921
$ws.Range("T2").Value = 0.003918177616359209
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9750066666666667
$ws.Range("H3").Value = 2.92502
$ws.Range("I3").Value = 0.03314938328005726
$ws.Range("J3").Value = 0.03314938328005725
$ws.Range("M3").Value = 0.7217316666666668
$ws.Range("O3").Value = 0.01768353618551768
$ws.Range("P3").Value = 0.01768353618551768
$ws.Range("Q3").Value = 0.7036931865444446
$ws.Range("R3").Value = 6.333238678900001
$ws.Range("S3").Value = 0.0005861983187604873
$ws.Range("T3").Value = 0.0005861983187604871
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9750066666666667
$ws.Range("H4").Value = 2.92502
$ws.Range("I4").Value = 0.03314938328005726
$ws.Range("J4").Value = 0.03314938328005725
$ws.Range("M4").Value = 4.718667333333333
$ws.Range("N4").Value = 14.156002
$ws.Range("O4").Value = 0.1156146091272429
$ws.Range("P4").Value = 0.1156146091272429
$ws.Range("Q4").Value = 4.600732107782222
$ws.Range("R4").Value = 41.40658897004
$ws.Range("S4").Value = 0.00383255299073298
$ws.Range("T4").Value = 0.003832552990732979
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9750066666666667
$ws.Range("H5").Value = 2.92502
$ws.Range("I5").Value = 0.03314938328005726
$ws.Range("J5").Value = 0.03314938328005725
$ws.Range("M5").Value = 27.97761966666667
$ws.Range("N5").Value = 83.93285900000001
$ws.Range("O5").Value = 0.6854947241613126
$ws.Range("P5").Value = 0.6854947241613125
$ws.Range("Q5").Value = 27.27836569246445
$ws.Range("R5").Value = 245.50529123218
$ws.Range("S5").Value = 0.02272372734768048
$ws.Range("T5").Value = 0.02272372734768047
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9750066666666667
$ws.Range("H6").Value = 2.92502
$ws.Range("I6").Value = 0.03314938328005726
$ws.Range("J6").Value = 0.03314938328005725
$ws.Range("M6").Value = 2.571656
$ws.Range("N6").Value = 7.714968
$ws.Range("O6").Value = 0.0630095283787885
$ws.Range("P6").Value = 0.06300952837878848
$ws.Range("Q6").Value = 2.507381744373333
$ws.Range("R6").Value = 22.56643569936
$ws.Range("S6").Value = 0.002088727006524105
$ws.Range("T6").Value = 0.002088727006524104
$ws.Range("I7").Value = 0.9153383232957121
$ws.Range("J7").Value = 0.9153383232957121
$ws.Range("M7").Value = 4.824089
$ws.Range("N7").Value = 14.472267
$ws.Range("O7").Value = 0.1181976021471384
$ws.Range("P7").Value = 0.1181976021471384
$ws.Range("Q7").Value = 129.8760553016297
$ws.Range("R7").Value = 1168.884497714667
$ws.Range("S7").Value = 0.1081907949669353
$ws.Range("T7").Value = 0.1081907949669353
$ws.Range("I8").Value = 0.9153383232957121
$ws.Range("J8").Value = 0.9153383232957121
$ws.Range("M8").Value = 0.7217316666666668
$ws.Range("O8").Value = 0.01768353618551768
$ws.Range("P8").Value = 0.01768353618551768
$ws.Range("Q8").Value = 19.43074886324389
$ws.Range("S8").Value = 0.01618641836199081
$ws.Range("T8").Value = 0.0161864183619908
$ws.Range("I9").Value = 0.9153383232957121
$ws.Range("J9").Value = 0.9153383232957121
$ws.Range("M9").Value = 4.718667333333333
$ws.Range("N9").Value = 14.156002
$ws.Range("O9").Value = 0.1156146091272429
$ws.Range("P9").Value = 0.1156146091272429
$ws.Range("Q9").Value = 127.0378509878224
$ws.Range("R9").Value = 1143.340658890402
$ws.Range("S9").Value = 0.1058264824670196
$ws.Range("T9").Value = 0.1058264824670196
$ws.Range("I10").Value = 0.9153383232957121
$ws.Range("J10").Value = 0.9153383232957121
$ws.Range("M10").Value = 27.97761966666667
$ws.Range("N10").Value = 83.93285900000001
$ws.Range("O10").Value = 0.6854947241613126
$ws.Range("P10").Value = 0.6854947241613125
$ws.Range("Q10").Value = 753.22467703974
$ws.Range("R10").Value = 6779.022093357659
$ws.Range("S10").Value = 0.6274595914418726
$ws.Range("T10").Value = 0.6274595914418725
$ws.Range("I11").Value = 0.9153383232957121
$ws.Range("J11").Value = 0.9153383232957121
$ws.Range("M11").Value = 2.571656
$ws.Range("N11").Value = 7.714968
$ws.Range("O11").Value = 0.0630095283787885
$ws.Range("P11").Value = 0.06300952837878848
$ws.Range("Q11").Value = 69.23515235161867
$ws.Range("R11").Value = 623.116371164568
$ws.Range("S11").Value = 0.05767503605789385
$ws.Range("T11").Value = 0.05767503605789383
$ws.Range("G12").Value = 1.515106
$ws.Range("H12").Value = 4.545318
$ws.Range("I12").Value = 0.05151229342423071
$ws.Range("J12").Value = 0.0515122934242307
$ws.Range("M12").Value = 4.824089
$ws.Range("N12").Value = 14.472267
$ws.Range("O12").Value = 0.1181976021471384
$ws.Range("P12").Value = 0.1181976021471384
$ws.Range("Q12").Value = 7.309006188434
$ws.Range("R12").Value = 65.78105569590601
$ws.Range("S12").Value = 0.006088629563843875
$ws.Range("T12").Value = 0.006088629563843874
$ws.Range("G13").Value = 1.515106
$ws.Range("H13").Value = 4.545318
$ws.Range("I13").Value = 0.05151229342423071
$ws.Range("J13").Value = 0.0515122934242307
$ws.Range("M13").Value = 0.7217316666666668
$ws.Range("O13").Value = 0.01768353618551768
$ws.Range("P13").Value = 0.01768353618551768
$ws.Range("Q13").Value = 1.093499978556667
$ws.Range("R13").Value = 9.841499807010001
$ws.Range("S13").Value = 0.0009109195047663881
$ws.Range("T13").Value = 0.0009109195047663879
$ws.Range("G14").Value = 1.515106
$ws.Range("H14").Value = 4.545318
$ws.Range("I14").Value = 0.05151229342423071
$ws.Range("J14").Value = 0.0515122934242307
$ws.Range("M14").Value = 4.718667333333333
$ws.Range("N14").Value = 14.156002
$ws.Range("O14").Value = 0.1156146091272429
$ws.Range("P14").Value = 0.1156146091272429
$ws.Range("Q14").Value = 7.149281188737334
$ws.Range("R14").Value = 64.343530698636
$ws.Range("S14").Value = 0.005955573669490276
$ws.Range("T14").Value = 0.005955573669490275
$ws.Range("G15").Value = 1.515106
$ws.Range("H15").Value = 4.545318
$ws.Range("I15").Value = 0.05151229342423071
$ws.Range("J15").Value = 0.0515122934242307
$ws.Range("M15").Value = 27.97761966666667
$ws.Range("N15").Value = 83.93285900000001
$ws.Range("O15").Value = 0.6854947241613126
$ws.Range("P15").Value = 0.6854947241613125
$ws.Range("Q15").Value = 42.38905942268467
$ws.Range("R15").Value = 381.501534804162
$ws.Range("S15").Value = 0.03531140537175963
$ws.Range("T15").Value = 0.03531140537175961
$ws.Range("G16").Value = 1.515106
$ws.Range("H16").Value = 4.545318
$ws.Range("I16").Value = 0.05151229342423071
$ws.Range("J16").Value = 0.0515122934242307
$ws.Range("M16").Value = 2.571656
$ws.Range("N16").Value = 7.714968
$ws.Range("O16").Value = 0.0630095283787885
$ws.Range("P16").Value = 0.06300952837878848
$ws.Range("Q16").Value = 3.896331435536
$ws.Range("R16").Value = 35.066982919824
$ws.Range("S16").Value = 0.003245765314370545
$ws.Range("T16").Value = 0.003245765314370544
